$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (closest reachable via ColumnWidth 1/6-character grid)
$ws.Range("A1").ColumnWidth = 14.833333333333334
$ws.Range("B1").ColumnWidth = 15.666666666666666

# Updated cell values
$ws.Range("A1").Value = -0.09611372762975634
$ws.Range("B1").Value = 0.095632035851991759
$ws.Range("A2").Value = -0.07093907640559749
$ws.Range("B2").Value = 0.069638325377807675
$ws.Range("A3").Value = -0.019924743359849728
$ws.Range("B3").Value = 0.019564932844961191
$ws.Range("A4").Value = -0.011564932965711705
$ws.Range("B4").Value = 0.011243167052221281
$ws.Range("A5").Value = -0.008243167103029414
$ws.Range("B5").Value = 0.0071522550532483464
$ws.Range("A6").Value = -0.0095491312260396199
$ws.Range("B6").Value = 0.0093270282800723692
$ws.Range("A7").Value = 0.00067297157617796088
$ws.Range("B7").Value = -0.00071607393566530675
$ws.Range("A8").Value = 0.010716073792936598
$ws.Range("B8").Value = -0.01077475320326915
$ws.Range("A9").Value = 0.012774753177312803
$ws.Range("B9").Value = -0.012818039667339853
$ws.Range("A10").Value = 0.014818039647055414
$ws.Range("B10").Value = -0.014817746708713386
$ws.Range("A11").Value = 0.017817746674862356
$ws.Range("B11").Value = -0.017820914949085775
$ws.Range("A12").Value = 0.021320914909615407
$ws.Range("B12").Value = -0.021364275260494558
$ws.Range("A13").Value = -0.0058443581911626907
$ws.Range("B13").Value = 0.0058372763082186196
$ws.Range("A14").Value = 0.0021627236009438278
$ws.Range("B14").Value = -0.0021626248503290313
$ws.Range("A15").Value = 0.0031626248579179617
$ws.Range("B15").Value = -0.0031637921910423472
$ws.Range("A16").Value = -0.0060336474915945182
$ws.Range("B16").Value = 0.0060029500925953272
$ws.Range("A17").Value = -0.0040029500986520361
$ws.Range("B17").Value = 0.0039999999657345242
$ws.Range("A18").Value = -0.018692959531389164
$ws.Range("B18").Value = 0.018579526026609727
$ws.Range("A19").Value = -0.012090786531534281
$ws.Range("B19").Value = 0.012015662772189462
$ws.Range("A20").Value = -0.0080156628344241199
$ws.Range("B20").Value = 0.0080056471551515074
$ws.Range("A21").Value = -0.0040056472180642899
$ws.Range("B21").Value = 0.0039999999365107897
$ws.Range("A22").Value = -0.045713582082461457
$ws.Range("B22").Value = 0.045500526265692898
$ws.Range("A23").Value = -0.040500526348999699
$ws.Range("B23").Value = 0.040099067745019923
$ws.Range("A24").Value = -0.020099068045637658
$ws.Range("B24").Value = 0.019999999695041737
$ws.Range("A25").Value = -0.079827400245683577
$ws.Range("B25").Value = 0.079761923346470809
$ws.Range("A26").Value = -0.077261923390139486
$ws.Range("B26").Value = 0.077177473117451711
$ws.Range("A27").Value = -0.074677473165066566
$ws.Range("B27").Value = 0.074174503199327724
$ws.Range("A28").Value = -0.072174503261246414
$ws.Range("B28").Value = 0.071838111255617676
$ws.Range("A29").Value = -0.064838111399138754
$ws.Range("B29").Value = 0.06474816007787787
$ws.Range("A30").Value = -0.0047481609679111969
$ws.Range("B30").Value = 0.0047388581363594362
$ws.Range("A31").Value = -0.014020893497074027
$ws.Range("B31").Value = 0.014000699274763662
$ws.Range("A32").Value = -0.0040006994690671149
$ws.Range("B32").Value = 0.0039999998892401578
